$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new cell values on row 5 (K5, L5)
$ws.Range("K5").Value = "Länge"
$ws.Range("L5").Value = "Typ Strom"

# Update the selection / active cell to L4:P4
$ws.Range("L4:P4").Select()

# Note: the author's saved window geometry (bookViews/workbookView
# windowWidth/windowHeight) and the xr:revisionPtr documentId GUID are
# Excel-generated, session-local artifacts of the save process. They are
# not exposed anywhere on the Application/Workbook/Window COM surface
# (confirmed: assigning ActiveWindow.Width/Height is accepted but produces
# no model mutation/diff at all), so there is no COM-interop call that can
# reproduce them here.
